$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H2: "1. manual / 2. auto" -> "manual", center-aligned + wrap (new style)
$ws.Range("H2").Value = "manual"
$ws.Range("H2").HorizontalAlignment = -4108  # xlCenter
$ws.Range("H2").WrapText = $true

# Duplicate row 2's formatting onto row 3 (reuses existing style indices)
$ws.Range("A2:H2").Copy()
$ws.Range("A3:H3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New "Sample B" test case row
$ws.Range("A3").Value = "Sample B"
$ws.Range("B3").Value = "This is a sample Test case"
$ws.Range("C3").Value = "medium"
$ws.Range("D3").Value = "None"
$ws.Range("E3").Value = "None"
$ws.Range("F3").Value = "1. Go to Google`n2. Display the object"
$ws.Range("G3").Value = "As Expected"
$ws.Range("H3").Value = "manual"

# Match row 2's wrapped-text row height
$ws.Rows.Item(3).RowHeight = 30

# Update active selection
$ws.Range("E10").Select()
